$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 30548.584
$ws.Range("J93").Value = 30548.584
$ws.Range("L93").Value = 30548.584
$ws.Range("N93").Value = -35540.584
$ws.Range("H95").Value = 34996
$ws.Range("J95").Value = 34996
$ws.Range("L95").Value = 34996
$ws.Range("N95").Value = -40488

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H95").Value = 40354
$ws.Range("J95").Value = 40354
$ws.Range("L95").Value = 40354
$ws.Range("N95").Value = -45846
$ws.Range("H104").Value = 42205
$ws.Range("J104").Value = 42205
$ws.Range("L104").Value = 42205
$ws.Range("N104").Value = -49193
$ws.Range("H105").Value = 47996
$ws.Range("J105").Value = 47996
$ws.Range("L105").Value = 47996
$ws.Range("N105").Value = -54984
$ws.Range("H106").Value = 47448.668
$ws.Range("J106").Value = 47448.668
$ws.Range("L106").Value = 47448.668
$ws.Range("N106").Value = -49972.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 41331.668
$ws.Range("J103").Value = 41331.668
$ws.Range("L103").Value = 41331.668
$ws.Range("N103").Value = -43675.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 43661
$ws.Range("J43").Value = 43661
$ws.Range("L43").Value = 43661
$ws.Range("N43").Value = -44029
$ws.Range("H92").Value = 38119.75
$ws.Range("J92").Value = 39851.285
$ws.Range("L92").Value = 39851.285
$ws.Range("N92").Value = -44843.285
$ws.Range("H96").Value = 73310
$ws.Range("J96").Value = 73310
$ws.Range("L96").Value = 73310
$ws.Range("N96").Value = -78802
$ws.Range("H101").Value = 43661
$ws.Range("J101").Value = 43661
$ws.Range("L101").Value = 43661
$ws.Range("N101").Value = -50151
$ws.Range("H106").Value = 32726.5
$ws.Range("J106").Value = 32726.5
$ws.Range("L106").Value = 32726.5
$ws.Range("N106").Value = -35250.5
$ws.Range("H124").Value = 45318
$ws.Range("J124").Value = 45318
$ws.Range("L124").Value = 45318
$ws.Range("N124").Value = -50228
$ws.Range("H125").Value = 49318
$ws.Range("J125").Value = 49318
$ws.Range("L125").Value = 49318
$ws.Range("N125").Value = -54238
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 35006.8
$ws.Range("J98").Value = 35006.8
$ws.Range("L98").Value = 35006.8
$ws.Range("N98").Value = -40996.8
$ws.Range("H100").Value = 38999
$ws.Range("J100").Value = 38999
$ws.Range("L100").Value = 38999
$ws.Range("N100").Value = -41163
$ws.Range("H104").Value = 35250.75
$ws.Range("J104").Value = 35250.75
$ws.Range("L104").Value = 35250.75
$ws.Range("N104").Value = -42238.75
$ws.Range("H105").Value = 44664.168
$ws.Range("J105").Value = 44664.168
$ws.Range("L105").Value = 44664.168
$ws.Range("N105").Value = -51652.168
$ws.Range("H118").Value = 38302
$ws.Range("J118").Value = 38302
$ws.Range("L118").Value = 38302
$ws.Range("N118").Value = -41616
$ws.Range("H120").Value = 39317
$ws.Range("J120").Value = 39317
$ws.Range("L120").Value = 39317
$ws.Range("N120").Value = -48993
$ws.Range("H125").Value = 40996
$ws.Range("J125").Value = 40996
$ws.Range("L125").Value = 40996
$ws.Range("N125").Value = -45916
$ws.Range("H130").Value = 48108.223
$ws.Range("J130").Value = 48108.223
$ws.Range("L130").Value = 48108.223
$ws.Range("N130").Value = -58148.223
$ws.Range("H131").Value = 38986
$ws.Range("J131").Value = 38986
$ws.Range("L131").Value = 38986
$ws.Range("N131").Value = -49066

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H95").Value = 35333
$ws.Range("J95").Value = 35333
$ws.Range("L95").Value = 35333
$ws.Range("N95").Value = -40825
$ws.Range("H98").Value = 46100
$ws.Range("J98").Value = 46100
$ws.Range("L98").Value = 46100
$ws.Range("N98").Value = -52090
$ws.Range("H103").Value = 48040.668
$ws.Range("J103").Value = 48040.668
$ws.Range("L103").Value = 48040.668
$ws.Range("N103").Value = -50384.668
$ws.Range("H105").Value = 46896
$ws.Range("J105").Value = 46896
$ws.Range("L105").Value = 46896
$ws.Range("N105").Value = -53884
$ws.Range("H117").Value = 45388
$ws.Range("J117").Value = 45388
$ws.Range("L117").Value = 45388
$ws.Range("N117").Value = -54566
$ws.Range("H123").Value = 32872.668
$ws.Range("J123").Value = 32872.668
$ws.Range("L123").Value = 32872.668
$ws.Range("N123").Value = -42672.668
$ws.Range("H127").Value = 41501.855
$ws.Range("J127").Value = 46643.832
$ws.Range("L127").Value = 46643.832
$ws.Range("N127").Value = -56563.832
$ws.Range("H131").Value = 43318
$ws.Range("J131").Value = 43318
$ws.Range("L131").Value = 43318
$ws.Range("N131").Value = -53398

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 28906.5
$ws.Range("J27").Value = 28906.5
$ws.Range("L27").Value = 28906.5
$ws.Range("N27").Value = -29044.5
$ws.Range("H95").Value = 35596.8
$ws.Range("J95").Value = 35596.8
$ws.Range("L95").Value = 35596.8
$ws.Range("N95").Value = -41088.8
$ws.Range("H97").Value = 36319
$ws.Range("J97").Value = 36319
$ws.Range("L97").Value = 36319
$ws.Range("N97").Value = -38301
$ws.Range("H98").Value = 36792.4
$ws.Range("J98").Value = 36792.4
$ws.Range("L98").Value = 36792.4
$ws.Range("N98").Value = -42782.4
$ws.Range("H103").Value = 47913.2
$ws.Range("J103").Value = 47913.2
$ws.Range("L103").Value = 47913.2
$ws.Range("N103").Value = -50257.2
$ws.Range("H104").Value = 44185
$ws.Range("J104").Value = 44185
$ws.Range("L104").Value = 44185
$ws.Range("N104").Value = -51173
$ws.Range("H115").Value = 37369
$ws.Range("J115").Value = 37369
$ws.Range("L115").Value = 37369
$ws.Range("N115").Value = -40503
$ws.Range("H118").Value = 37694.668
$ws.Range("J118").Value = 37694.668
$ws.Range("L118").Value = 37694.668
$ws.Range("N118").Value = -41008.668
$ws.Range("H129").Value = 27221
$ws.Range("J129").Value = 27221
$ws.Range("L129").Value = 27221
$ws.Range("N129").Value = -37221
